# 0709: modify description errors
#
# Slide 1 holds a notice textbox ("TextBox 2") nested two levels deep
# inside grouped shapes ("그룹 7" > "그룹 9" siblings). PowerPoint's
# GroupItems collection flattens that nesting, so the textbox is item 6
# of the top-level group's GroupItems.
#
# Paragraph 1 of that textbox announces an ACL acceptance and needs two
# text corrections (same-length replacements, so visually nothing
# reflows):
#   - "ACL 2024"                         -> "ACL 2025"
#   - "(BK Plus Computer Science IF=3)"  -> "(BK Plus Computer Science IF=4)"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$topGroup = $s.Shapes.Item(1)
$noticeBox = $topGroup.GroupItems.Item(6)

# The textbox auto-fits its height to the text ("Resize shape to fit
# text"). Re-assigning run text makes the host re-run that auto-fit
# pass, which can nudge the cached shape height by a hair even though
# these two same-length substitutions don't actually change the
# wrapped line count. Snapshot the exact height up front and restore it
# after the edit so the shape geometry is untouched, matching the
# source diff (only the two text runs changed).
$origHeightPt = $noticeBox.Height

$tr = $noticeBox.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

$para1.Runs(1, 1).Text = "ACL 2025"
$para1.Runs(5, 1).Text = "(BK Plus Computer Science IF=4)"

# Restore the original height precisely. Toggling AutoSize off/on
# around the Height write keeps the auto-fit pass from re-stamping its
# own (slightly different) recalculated height over ours; the Height
# assignment must come after AutoSize is switched back on, since that
# switch itself re-triggers a fit pass. The host stores shape geometry
# in EMUs but the Height property is in points (1 pt = 12700 EMU) and
# round-trips through single precision, so naively writing back the pt
# value read above truncates down by one EMU. Nudging by half an EMU's
# worth of points compensates for that truncation and lands back on the
# exact original EMU value.
$noticeBox.TextFrame.AutoSize = 0
$noticeBox.TextFrame.AutoSize = 1
$noticeBox.Height = $origHeightPt + (0.5 / 12700.0)
